$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: split the text run at [splitStart, splitEnd) away from its
# neighbours by toggling Bold on/off (net no-op formatting change) so the
# engine keeps it as its own <w:r> instead of silently re-coalescing it with
# an adjacent run that happens to share identical formatting.
# ---------------------------------------------------------------------------
function Pin-Run($rng) {
    $rng.Bold = $true
    $rng.Bold = $false
}

# ===========================================================================
# Change 1: "http://server:port" -> "http://" + "<domain>"
#   (web page URL, first mention)
# ===========================================================================
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("http://server:port", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$matchStart = $rng.Start
$matchEnd = $rng.End

$httpLen = 7   # Len("http://")
$r2 = $d.Range($matchStart + $httpLen, $matchEnd)
$r2.Text = "<domain>"
$r2 = $d.Range($matchStart + $httpLen, $matchStart + $httpLen + 8)  # Len("<domain>")
Pin-Run $r2

Write-Host "Change 1 done"

# ===========================================================================
# Change 2: "http://" + "server:port" -> "http:" + "//" + "<domain>"
#   (service URL mention, "http://server:port/mazes/render")
#   The trailing "/mazes/render" run is untouched.
# ===========================================================================
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("http://server:port", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$matchStart = $rng.Start
$matchEnd = $rng.End   # matchEnd - matchStart == Len("http://server:port") == 19

# "server:port" (positions matchStart+7 .. matchEnd) becomes "<domain>"
$rDomain = $d.Range($matchStart + 7, $matchEnd)
$rDomain.Text = "<domain>"

# Re-derive the "<domain>" range (8 chars) and pin it off from the "//" run.
$rDomain = $d.Range($matchStart + 7, $matchStart + 7 + 8)
Pin-Run $rDomain

# Split "http://" into "http:" (5 chars) + "//" (2 chars); pin the "//" part
# off from the preceding "http:" run so it stays separate.
$rSlashes = $d.Range($matchStart + 5, $matchStart + 7)
Pin-Run $rSlashes

Write-Host "Change 2 done"

# ===========================================================================
# Change 3: "http://server:port/mazes/render" -> "http://" + "<domain>" +
#   "/mazes/render"   (the "POST <url>" example line)
# ===========================================================================
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("http://server:port/mazes/render", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$matchStart = $rng.Start
$matchEnd = $rng.End

# "server:port" -> "<domain>"
$rDomain = $d.Range($matchStart + 7, $matchStart + 7 + 11)  # Len("server:port") == 11
$rDomain.Text = "<domain>"

# Pin "http://" (7 chars) off from the preceding (identically formatted)
# " " run so it remains its own run.
$rHttp = $d.Range($matchStart, $matchStart + 7)
Pin-Run $rHttp

# Pin "<domain>" (8 chars) off from the preceding "http://" run.
$rDomain = $d.Range($matchStart + 7, $matchStart + 7 + 8)
Pin-Run $rDomain

Write-Host "Change 3 done"

# ===========================================================================
# Change 4: "As a bonus, you may also code up one, or more, maze solvers to
#   apply to your generated output." ->
#   "If time permits" + ", you may " + "want to" +
#   " code up one, or more, maze solvers to apply to your generated output."
# ===========================================================================
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("As a bonus, you may also code up one, or more, maze solvers to apply to your generated output.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$matchStart = $rng.Start

$newText = "If time permits, you may want to code up one, or more, maze solvers to apply to your generated output."
$rng.Text = $newText

$p1 = "If time permits"
$p2 = ", you may "
$p3 = "want to"
# p4 is the remainder of the text.

$b1 = $matchStart + $p1.Length
$b2 = $b1 + $p2.Length
$b3 = $b2 + $p3.Length

# Pin each of the first three pieces off from the run that follows it so the
# four segments end up as four distinct runs.
Pin-Run ($d.Range($matchStart, $b1))
Pin-Run ($d.Range($b1, $b2))
Pin-Run ($d.Range($b2, $b3))

Write-Host "Change 4 done"
